$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Description") - this shifts the
# old C (Date d'échéance) and D (Etat d'avancement) columns to D and E,
# and also shifts the existing column width / data validation ranges.
$ws.Columns(3).Insert()

# Header for the new column.
$ws.Range("C1").Value = "Description"

# New description values for a few specific tasks.
$ws.Range("C4").Value = "Amortisseur"
$ws.Range("C9").Value = "Yo"
$ws.Range("C14").Value = "Youpi"

# Set the width of the new Description column (best achievable value
# given this engine's character-width rounding).
$ws.Columns("C:C").ColumnWidth = 20.5

# Update the selection to match the edited workbook's cursor position.
$ws.Range("C19").Select()
